# CISPR_25.xlsx edit:
# Insert a new header row at the top of the "CISPR 25" sheet ("F, Mhz" /
# "Eeq max ref"), which pushes the existing 481 data rows down by one
# (rows 1-481 -> rows 2-482). Update the CISPR_25 defined name range and
# the active cell selection to reflect the new layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift all existing data down by one row to make room for the header.
$ws.Rows.Item(1).Insert()

# Fill in the new header row.
$ws.Range("A1").Value = "F, Mhz"
$ws.Range("B1").Value = "Eeq max ref"
$ws.Rows.Item(1).RowHeight = 30

# The CISPR_25 named range used to cover the data at $A$1:$B$481; the
# data now lives one row down, at $A$2:$B$482.
$n = $wb.Names.Item("CISPR_25")
$n.RefersTo = "='CISPR 25'!`$A`$2:`$B`$482"

# Update the selected cell shown when the sheet is opened.
$ws.Range("B2").Select()
